# Generate Report for Archive
#
# The "Status" value that used to read "Ready for handoff" is now
# "In Translation" everywhere it appears (Overview!E2/F2 - the per-language
# status columns - and the Status column (C2) on each per-language sheet).
# Because the new text is shorter, the report generator that produced this
# workbook re-sized ("auto-fit") the columns that held that value, shrinking
# them from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# New, narrower width those auto-fitted columns were resized to (target
# stored OOXML width is 13.4101845877511 "characters"; this runtime quantizes
# ColumnWidth to sixths of a character internally, so 12.5 is the input that
# rounds to the stored width nearest that target, 13.333333333333334).
$newStatusColWidth = 12.5

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet: column C ("Status") holds the status text
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet: column C ("Status") holds the status text
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
